# Nexial Interactive alpha - adds sendHtmlMail(...) and clickWithKeys(...) commands
# to the hidden #system reference sheet, and updates the "aws.ses" sample command.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Replace the sample "aws.ses" command shown in C2 (sendMail -> sendHtmlMail)
$ws.Range("C2").Value = "sendHtmlMail(profile,to,subject,body)"

# 2) Insert the new "clickWithKeys(locator,keys)" command into the alphabetically
#    sorted "web" command list (column V), pushing close() .. waitForTitle() down
#    by one row (V51..V118).
$newWebList = @(
    'clickWithKeys(locator,keys)',
    'close()',
    'closeAll()',
    'deselect(locator,text)',
    'deselectMulti(locator,array)',
    'dismissInvalidCert()',
    'dismissInvalidCertPopup()',
    'doubleClick(locator)',
    'doubleClickAndWait(locator,waitMs)',
    'doubleClickByLabel(label)',
    'doubleClickByLabelAndWait(label,waitMs)',
    'dragAndDrop(fromLocator,toLocator)',
    'editLocalStorage(key,value)',
    'executeScript(var,script)',
    'focus(locator)',
    'goBack()',
    'goBackAndWait()',
    'maximizeWindow()',
    'mouseOver(locator)',
    'open(url)',
    'openAndWait(url,waitMs)',
    'openHttpBasic(url,username,password)',
    'refresh()',
    'refreshAndWait()',
    'resizeWindow(width,height)',
    'saveAllWindowIds(var)',
    'saveAllWindowNames(var)',
    'saveAttribute(var,locator,attrName)',
    'saveCount(var,locator)',
    'saveDivsAsCsv(headers,rows,cells,nextPage,file)',
    'saveElement(var,locator)',
    'saveElements(var,locator)',
    'saveLocalStorage(var,key)',
    'saveLocation(var)',
    'savePageAs(var,sessionIdName,url)',
    'savePageAsFile(sessionIdName,url,file)',
    'saveTableAsCsv(locator,nextPageLocator,file)',
    'saveText(var,locator)',
    'saveTextArray(var,locator)',
    'saveTextSubstringAfter(var,locator,delim)',
    'saveTextSubstringBefore(var,locator,delim)',
    'saveTextSubstringBetween(var,locator,start,end)',
    'saveValue(var,locator)',
    'scrollLeft(locator,pixel)',
    'scrollRight(locator,pixel)',
    'scrollTo(locator)',
    'select(locator,text)',
    'selectFrame(locator)',
    'selectMulti(locator,array)',
    'selectMultiOptions(locator)',
    'selectText(locator)',
    'selectWindow(winId)',
    'selectWindowAndWait(winId,waitMs)',
    'selectWindowByIndex(index)',
    'selectWindowByIndexAndWait(index,waitMs)',
    'toggleSelections(locator)',
    'type(locator,value)',
    'typeKeys(locator,value)',
    'uncheckAll(locator)',
    'unselectAllText()',
    'upload(fieldLocator,file)',
    'verifyContainText(locator,text)',
    'verifyText(locator,text)',
    'wait(waitMs)',
    'waitForElementPresent(locator)',
    'waitForPopUp(winId,waitMs)',
    'waitForTextPresent(text)',
    'waitForTitle(text)'
)

for ($i = 0; $i -lt $newWebList.Length; $i++) {
    $row = 51 + $i
    $ws.Cells.Item($row, 22).Value = $newWebList[$i]
}

# 3) Extend the "web" defined name range to cover the new last row (118).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "web") {
        $n.RefersTo = "='#system'!`$V`$2:`$V`$118"
    }
}
